# Add a new test site row (row 3) to the generic_one_table "pages" sheet,
# cloning the formatting of the existing "vgchartz" row (row 2) and giving
# the new row a distinct nickname "vgchartz_2".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the cell formatting (styles/borders/fonts/alignment) of row 2 onto row 3.
$ws.Range("A2:G2").Copy()
$ws.Range("A3:G3").PasteSpecial(-4122)  # xlPasteFormats

# Populate row 3's values: same url / browser / html selectors as row 2,
# but a new nickname to distinguish this as a second test entry.
$ws.Range("A3").Value = $ws.Range("A2").Text
$ws.Range("B3").Value = $ws.Range("B2").Text
$ws.Range("C3").Value = "vgchartz_2"
$ws.Range("D3").Value = $ws.Range("D2").Text
$ws.Range("E3").Value = $ws.Range("E2").Text
$ws.Range("F3").Value = $ws.Range("F2").Text
$ws.Range("G3").Value = $ws.Range("G2").Text
